# Mise à jour de l'application
# Add a new attendance date column (AN) for 2025-10-09 (serial 45904),
# one column to the right of the existing last date column (AM).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New date header cell (AN1) -------------------------------------------------
$ws.Range("AN1").Value2 = 45904

# --- New attendance values for each player row (AN2:AN27) ----------------------
$values = @{
    2  = "P"
    3  = "P"
    4  = "P"
    5  = "RH"
    6  = "B"
    7  = "P"
    8  = "P"
    9  = "P"
    10 = "B"
    11 = "RH"
    12 = "P"
    13 = "RH"
    14 = "P"
    15 = "P"
    16 = "P"
    17 = "RH"
    18 = "P"
    19 = "P"
    20 = "P"
    21 = "P"
    22 = "P"
    23 = "P"
    24 = "P"
    25 = "P"
    26 = "P"
    27 = "RH"
}

foreach ($row in $values.Keys) {
    $ws.Range("AN$row").Value2 = $values[$row]
}

# --- Match the formatting of the preceding date column (AM) --------------------
$ws.Range("AM1").Copy()
$ws.Range("AN1").PasteSpecial(-4122)

$ws.Range("AM2:AM27").Copy()
$ws.Range("AN2:AN27").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Update the active selection as in the saved file --------------------------
$ws.Range("AP23").Select()

Write-Output "done"
